# Update the "想去人数" (F column) values on both the "展览" and
# "全部类型" worksheets to reflect the latest scrape output.
#
# row -> new value for column F
$updates = @{
    2  = 269
    3  = 287
    4  = 11035
    5  = 10214
    13 = 9583
    15 = 2437
    17 = 6
    18 = 83
    19 = 389
    20 = 10857
    21 = 10779
    26 = 16
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
